# This script updates the "想去人数" (F column) figures on the
# "展览" and "全部类型" worksheets to reflect newly generated output
# (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value  = 15059   # F2  15056 -> 15059
$ws1.Cells.Item(3, 6).Value  = 19198   # F3  19184 -> 19198
$ws1.Cells.Item(22, 6).Value = 8017    # F22 8010  -> 8017
$ws1.Cells.Item(27, 6).Value = 1250    # F27 1248  -> 1250
$ws1.Cells.Item(29, 6).Value = 6080    # F29 6076  -> 6080
$ws1.Cells.Item(35, 6).Value = 5482    # F35 5476  -> 5482
$ws1.Cells.Item(36, 6).Value = 934     # F36 894   -> 934
$ws1.Cells.Item(38, 6).Value = 30      # F38 29    -> 30
$ws1.Cells.Item(39, 6).Value = 51      # F39 50    -> 51

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value  = 15059   # F2  15056 -> 15059
$ws4.Cells.Item(3, 6).Value  = 19199   # F3  19184 -> 19199
$ws4.Cells.Item(23, 6).Value = 8017    # F23 8010  -> 8017
$ws4.Cells.Item(28, 6).Value = 1250    # F28 1248  -> 1250
$ws4.Cells.Item(32, 6).Value = 6080    # F32 6076  -> 6080
$ws4.Cells.Item(38, 6).Value = 5482    # F38 5476  -> 5482
$ws4.Cells.Item(39, 6).Value = 934     # F39 894   -> 934
$ws4.Cells.Item(41, 6).Value = 30      # F41 29    -> 30
$ws4.Cells.Item(42, 6).Value = 51      # F42 50    -> 51
